# Add a new "Event" column (column M) to the Card24 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Column M is currently empty (sheet only uses A:L), so pick up the header
# look-and-feel (bold, centered, bordered) from the last existing header
# cell (L1) and apply it to the new header cell M1.
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").Value = "Event"
